# Updates cryptos list values per the commit diff: new Price (D) / Volume(1h) (E)
# figures for every coin row, plus two coin re-rankings (rows 38-39 swap
# SuiNetwork/Stacks, rows 42-43 swap Bittensor/FirstDigitalUSD) that moved the
# Coin (B) and Link (C) text between rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain text in this sheet (e.g. "6.48", "59.076.25")
# rather than numbers. Excel.Value auto-converts plain decimal-looking text to a
# real number, which would corrupt values like "6.49" -> 6.49 (number) and drop
# the inline-string typing the source file uses. Force those through as text via
# a quote-prefix, then reset Style so no stray quote-prefix formatting lingers.
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = 'Normal'
}

$ws.Range('D2').Value = '59.163.98'
$ws.Range('E2').Value = '  +2.49%  '
$ws.Range('D3').Value = '2.596.57'
$ws.Range('E3').Value = '  +1.78%  '
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue 'D5' '520.15'
$ws.Range('E5').Value = '  +0.40%  '
Set-TextValue 'D6' '139.40'
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('E8').Value = '  +1.06%  '
$ws.Range('D9').Value = '2.624.76'
$ws.Range('E9').Value = '  +2.80%  '
Set-TextValue 'D10' '6.49'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('E11').Value = '  +2.46%  '
$ws.Range('E12').Value = '  +2.58%  '
$ws.Range('E13').Value = '  +1.83%  '
$ws.Range('D14').Value = '3.063.55'
$ws.Range('E14').Value = '  +2.13%  '
$ws.Range('D15').Value = '59.095.70'
$ws.Range('E15').Value = '  +2.38%  '
Set-TextValue 'D16' '20.40'
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('D17').Value = '2.622.15'
$ws.Range('E17').Value = '  +3.50%  '
$ws.Range('E18').Value = '  +0.55%  '
Set-TextValue 'D19' '339.30'
$ws.Range('E19').Value = '  +1.89%  '
Set-TextValue 'D20' '4.32'
$ws.Range('E20').Value = '  +1.08%  '
Set-TextValue 'D21' '10.18'
$ws.Range('E21').Value = '  +0.89%  '
Set-TextValue 'D22' '6.50'
$ws.Range('E22').Value = '  +6.37%  '
$ws.Range('E23').Value = '  -0.21%  '
Set-TextValue 'D24' '66.38'
$ws.Range('E24').Value = '  +2.31%  '
Set-TextValue 'D25' '0.167'
$ws.Range('E25').Value = '  +2.19%  '
$ws.Range('E26').Value = '  +0.98%  '
Set-TextValue 'D27' '0.996'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('E28').Value = '  +1.49%  '
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').Value = '0.0₃0725'
$ws.Range('E30').Value = '  -3.45%  '
$ws.Range('E31').Value = '  -3.51%  '
Set-TextValue 'D32' '18.86'
$ws.Range('E32').Value = '  +2.49%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  +0.13%  '
Set-TextValue 'D35' '3.99'
$ws.Range('E35').Value = '  +0.82%  '
$ws.Range('E36').Value = '  +0.43%  '
Set-TextValue 'D37' '36.34'
$ws.Range('E37').Value = '  +1.86%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D38' '1.46'
$ws.Range('E38').Value = '  +4.39%  '
$ws.Range('B39').Value = 'SuiNetwork'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue 'D39' '0.830'
$ws.Range('E39').Value = '  -0.24%  '
Set-TextValue 'D40' '0.830'
$ws.Range('E40').Value = '  +1.65%  '
Set-TextValue 'D41' '3.54'
$ws.Range('E41').Value = '  +2.59%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D42' '0.995'
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D43' '275.39'
$ws.Range('E43').Value = '  +6.00%  '
$ws.Range('E44').Value = '  +0.72%  '
Set-TextValue 'D45' '0.594'
$ws.Range('E45').Value = '  +2.87%  '
$ws.Range('E46').Value = '  -0.29%  '
Set-TextValue 'D47' '0.0521'
$ws.Range('E47').Value = '  +0.68%  '
Set-TextValue 'D48' '18.51'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D49').Value = '1.986.55'
$ws.Range('E49').Value = '  +0.26%  '
Set-TextValue 'D50' '4.60'
$ws.Range('E50').Value = '  +2.06%  '
$ws.Range('E51').Value = '  -0.48%  '
